$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph 9: "9.  After entering new record and saving, returns to Main
#    Activity but keyboard " -> split into three runs by inserting
#    "editing a record or " right after "9.  After ".
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("9.  After ")
$rng.Collapse(0)                      # wdCollapseEnd
$insStart = $rng.Start

$ins = $rng.Duplicate
$ins.InsertAfter("editing a record or ")
$insEnd = $ins.End

# Wedge temporary bookmarks at both boundaries of the freshly-inserted text
# so the new text keeps its own run and is not merged back into its
# neighbours; then remove the temporary bookmarks (the run split survives).
$wedge1 = $d.Range($insStart, $insStart)
$d.Bookmarks.Add("zzzTempWedge1", $wedge1)
$wedge2 = $d.Range($insEnd, $insEnd)
$d.Bookmarks.Add("zzzTempWedge2", $wedge2)
$d.Bookmarks.Item("zzzTempWedge1").Delete()
$d.Bookmarks.Item("zzzTempWedge2").Delete()

# ---------------------------------------------------------------------------
# 2) " visible, obscuring " -> " visible, obscurin" + "g ", with the
#    "_GoBack" bookmark relocated to the split point (it previously sat
#    just before "DONE 12." later in the document; adding it again under
#    the same name moves it and removes the old occurrence).
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$null = $rng2.Find.Execute(" visible, obscuring ")
$splitPos = $rng2.Start + 18          # " visible, obscurin" is 18 chars long
$goBackRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)
